$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values (prices) are written as plain text, matching the
# original inline-string cell type, rather than being auto-converted to
# numbers by Excel's normal "typed value" interpretation.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.301.38"
$ws.Range("E2").Value = "  +0.92%  "

$ws.Range("D3").Value = "1.679.83"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").Value = "218.01"
$ws.Range("E5").Value = "  +0.57%  "

$ws.Range("D6").Value = "0.5393"
$ws.Range("E6").Value = "  +5.71%  "

$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("D8").Value = "0.2688"
$ws.Range("E8").Value = "  +1.22%  "

$ws.Range("D9").Value = "0.06475"
$ws.Range("E9").Value = "  +0.96%  "

$ws.Range("D10").Value = "21.98"
$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("D11").Value = "0.07538"
$ws.Range("E11").Value = "  +1.55%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "4.529"

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.673.01"
$ws.Range("E13").Value = "  -0.89%  "

$ws.Range("D14").Value = "0.5786"
$ws.Range("E14").Value = "  -1.00%  "

$ws.Range("D15").Value = "0.000008451"
$ws.Range("E15").Value = "  -1.63%  "

$ws.Range("E16").Value = "  +0.43%  "

$ws.Range("D17").Value = "26.312.17"
$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("D18").Value = "4.906"
$ws.Range("E18").Value = "  -0.97%  "

$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").Value = "10.88"

$ws.Range("D21").Value = "190.98"
$ws.Range("E21").Value = "  -1.09%  "

$ws.Range("D22").Value = "6.219"

$ws.Range("D23").Value = "1.007"
$ws.Range("E23").Value = "  +0.28%  "

$ws.Range("D24").Value = "146.22"
$ws.Range("E24").Value = "  +0.84%  "

$ws.Range("D25").Value = "0.1289"
$ws.Range("E25").Value = "  +7.80%  "

$ws.Range("D26").Value = "7.838"

$ws.Range("E27").Value = "  +0.36%  "

$ws.Range("D28").Value = "0.06501"
$ws.Range("E28").Value = "  +1.52%  "

$ws.Range("E29").Value = "  +3.93%  "

$ws.Range("D30").Value = "1.322"
$ws.Range("E30").Value = "  +0.27%  "

$ws.Range("D31").Value = "3.580"
$ws.Range("E31").Value = "  +0.70%  "

$ws.Range("E32").Value = "  +1.45%  "

$ws.Range("E33").Value = "  +0.96%  "

$ws.Range("D34").Value = "1.035"
$ws.Range("E34").Value = "  +1.41%  "

$ws.Range("D35").Value = "0.6162"
$ws.Range("E35").Value = "  +0.86%  "

$ws.Range("D36").Value = "2.400"
$ws.Range("E36").Value = "  +1.36%  "

$ws.Range("D37").Value = "2.722"
$ws.Range("E37").Value = "  +0.49%  "

$ws.Range("D38").Value = "6.245"
$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("D39").Value = "1.111.45"
$ws.Range("E39").Value = "  +1.95%  "

$ws.Range("E40").Value = "  +1.04%  "

$ws.Range("E41").Value = "  +0.86%  "

$ws.Range("E42").Value = "  +0.66%  "

$ws.Range("E43").Value = "  -0.23%  "

$ws.Range("D44").Value = "1.829.52"
$ws.Range("E44").Value = "  +0.71%  "

$ws.Range("D45").Value = "0.00000000110"
$ws.Range("E45").Value = "  -2.62%  "

$ws.Range("D46").Value = "57.19"
$ws.Range("E46").Value = "  +1.28%  "

$ws.Range("D47").Value = "8.139"
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("D48").Value = "1.000"
$ws.Range("E48").Value = "  -0.55%  "

$ws.Range("D49").Value = "0.05274"
$ws.Range("E49").Value = "  +0.73%  "

$ws.Range("D50").Value = "0.4290"
$ws.Range("E50").Value = "  +0.13%  "

$ws.Range("E51").Value = "  +0.69%  "

# Reset column D styling back to the default (no explicit number format)
# now that the text values are safely stored, so we don't leave a stray
# "@" text-format style applied to the cells.
$ws.Range("D2:D51").Style = "Normal"

Write-Host "Applied cryptos list update"